$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": insert a new client row for
# "MIM CONSTRUFERRETERIA E IMPORTADORA SAS" right before the
# existing "MOROCHO BACUILIMA HILDA INES" row (row 12), pushing the
# rest of the table (and the trailing totals row) down by one.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(12).Insert()

$ws1.Range("A12").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Range("B12").Value = "MIM CONSTRUFERRETERIA E IMPORTADORA SAS"
$ws1.Range("C12:R12").Value = 0

# The trailing "X de 13" summary row is now row 16; bump the counts
# to reflect the extra row ("X de 14").
$ws1.Range("C16").Value = "0 de 14"
$ws1.Range("D16").Value = "0 de 14"
$ws1.Range("E16").Value = "0 de 14"
$ws1.Range("F16").Value = "0 de 14"
$ws1.Range("G16").Value = "0 de 14"
$ws1.Range("H16").Value = "3 de 14"
$ws1.Range("I16").Value = "1 de 14"
$ws1.Range("J16").Value = "0 de 14"
$ws1.Range("K16").Value = "0 de 14"
$ws1.Range("L16").Value = "0 de 14"
$ws1.Range("M16").Value = "2 de 14"
$ws1.Range("N16").Value = "0 de 14"
$ws1.Range("O16").Value = "0 de 14"
$ws1.Range("P16").Value = "0 de 14"
$ws1.Range("Q16").Value = "0 de 14"
$ws1.Range("R16").Value = "0 de 14"

# Column B ("CLIENTE") widens from 36 to 41 characters to fit the
# longer new client name.
$ws1.Columns.Item(2).ColumnWidth = 40.1

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL": mirror the same row insertion.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(12).Insert()

$ws2.Range("A12").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Range("B12").Value = "MIM CONSTRUFERRETERIA E IMPORTADORA SAS"
$ws2.Range("C12:G12").Value = 0

$ws2.Columns.Item(2).ColumnWidth = 40.1
